# Insert a new data row at row 356 (pushing existing rows 356-443 down to
# 357-444) and populate it with the new weekly price record for Cilantro.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 356..443 down by one row, growing the used range to row 444.
$ws.Rows.Item(356).Insert()

# Populate the newly inserted row 356 with the new record.
$ws.Range("A356").Value = 3
$ws.Range("B356").Value = "Femacal de La Calera"
$ws.Range("C356").Value = "Coquimbo"
$ws.Range("D356").Value = 44855
$ws.Range("E356").Value = 5
$ws.Range("F356").Value = 100112040
$ws.Range("G356").Value = "Cilantro"
$ws.Range("H356").Value = "Sin especificar"
$ws.Range("I356").Value = "Primera"
$ws.Range("J356").Value = 280
$ws.Range("K356").Value = 2500
$ws.Range("L356").Value = 2800
$ws.Range("M356").Value = 2629
$ws.Range("N356").Value = "$/docena de atados (3 kilos)"
$ws.Range("O356").Value = "Provincia de Quillota"
$ws.Range("P356").Value = 876
$ws.Range("Q356").Value = 3
$ws.Range("R356").Value = "Hortaliza"
